$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weekly_detail")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://sampling-458.netlify.app/notes/cn04-srs", "", "", "https://sampling-458.netlify.app/notes/cn04-srs")
